$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 9790
$ws.Range("F3").Value = 415
$ws.Range("F4").Value = 2518
$ws.Range("F9").Value = 715
$ws.Range("F13").Value = 3026
$ws.Range("F14").Value = 2293
$ws.Range("F16").Value = 1991
$ws.Range("F17").Value = 243
$ws.Range("F21").Value = 321
$ws.Range("F22").Value = 29
$ws.Range("F23").Value = 200
$ws.Range("F25").Value = 38
$ws.Range("F26").Value = 353
$ws.Range("F27").Value = 73
$ws.Range("F28").Value = 329
$ws.Range("F29").Value = 538
$ws.Range("F30").Value = 40
$ws.Range("F31").Value = 180
$ws.Range("F32").Value = 1555
$ws.Range("F33").Value = 222
$ws.Range("F34").Value = 1557
$ws.Range("F35").Value = 72
$ws.Range("F36").Value = 375
$ws.Range("F38").Value = 403
$ws.Range("F39").Value = 830
$ws.Range("F41").Value = 326

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 6

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9790
$ws.Range("F3").Value = 415
$ws.Range("F4").Value = 2518
$ws.Range("F11").Value = 715
$ws.Range("F15").Value = 3026
$ws.Range("F16").Value = 2293
$ws.Range("F18").Value = 1991
$ws.Range("F19").Value = 243
$ws.Range("F23").Value = 322
$ws.Range("F24").Value = 29
$ws.Range("F25").Value = 200
$ws.Range("F27").Value = 38
$ws.Range("F28").Value = 353
$ws.Range("F29").Value = 73
$ws.Range("F30").Value = 329
$ws.Range("F31").Value = 538
$ws.Range("F35").Value = 40
$ws.Range("F36").Value = 180
$ws.Range("F37").Value = 1555
$ws.Range("F39").Value = 222
$ws.Range("F40").Value = 1557
$ws.Range("F41").Value = 72
$ws.Range("F42").Value = 6
$ws.Range("F43").Value = 375
$ws.Range("F45").Value = 403
$ws.Range("F46").Value = 830
$ws.Range("F48").Value = 326
